$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> values for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$data = @{
  2  = @{ E=3; G=12.72626133333333; H=38.178784; I=0.08664518826379554; J=0.08664518826379553; K=3; M=43.027021; N=129.081063; O=0.7775798784610731; P=0.7775798784610732; Q=547.5731136408214; R=4928.158022767392; S=0.06737355495939894; T=0.06737355495939894 }
  3  = @{ E=3; G=12.72626133333333; H=38.178784; I=0.08664518826379554; J=0.08664518826379553; K=3; M=2.947472666666667; N=8.842418; O=0.05326642153343597; P=0.05326642153343598; Q=37.51030742885689; R=337.592766859712; S=0.004615279121903252; T=0.004615279121903252 }
  4  = @{ E=3; G=12.72626133333333; H=38.178784; I=0.08664518826379554; J=0.08664518826379553; K=3; M=9.360041333333333; N=28.080124; O=0.1691537000054908; P=0.1691537000054908; Q=119.1183320988018; R=1072.064988889216; S=0.01465635418249335; T=0.01465635418249334 }
  5  = @{ E=3; G=132.5150096666667; H=397.545029; I=0.9022121783931373; J=0.9022121783931372; K=3; M=43.027021; N=129.081063; O=0.7775798784610731; P=0.7775798784610732; Q=5701.726103742869; R=51315.53493368583; S=0.7015420360210357; T=0.7015420360210357 }
  6  = @{ E=3; G=132.5150096666667; H=397.545029; I=0.9022121783931373; J=0.9022121783931372; K=3; M=2.947472666666667; N=8.842418; O=0.05326642153343597; P=0.05326642153343598; Q=390.5843689155691; R=3515.259320240122; S=0.04805761420688838; T=0.04805761420688838 }
  7  = @{ E=3; G=132.5150096666667; H=397.545029; I=0.9022121783931373; J=0.9022121783931372; K=3; M=9.360041333333333; N=28.080124; O=0.1691537000054908; P=0.1691537000054908; Q=1240.345967767066; R=11163.1137099036; S=0.1526125281652131; T=0.1526125281652131 }
  8  = @{ E=3; G=1.636606333333333; H=4.909819; I=0.0111426333430672; J=0.01114263334306719; K=3; M=43.027021; N=129.081063; O=0.7775798784610731; P=0.7775798784610732; Q=70.41829507306632; R=633.764655657597; S=0.00866428748063849; T=0.00866428748063849 }
  9  = @{ E=3; G=1.636606333333333; H=4.909819; I=0.0111426333430672; J=0.01114263334306719; K=3; M=2.947472666666667; N=8.842418; O=0.05326642153343597; P=0.05326642153343598; Q=4.823852433593555; R=43.414671902342; S=0.000593528204644336; T=0.000593528204644336 }
  10 = @{ E=3; G=1.636606333333333; H=4.909819; I=0.0111426333430672; J=0.01114263334306719; K=3; M=9.360041333333333; N=28.080124; O=0.1691537000054908; P=0.1691537000054908; Q=15.31870292639511; R=137.868326337556; S=0.001884817657784368; T=0.001884817657784368 }
}

foreach ($row in $data.Keys) {
  $cols = $data[$row]
  foreach ($col in $cols.Keys) {
    $ws.Range("$col$row").Value = $cols[$col]
  }
}
